$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for column B (rows 2-158), replacing the placeholder "1" values
# with the actual computed scores.
$bValues = @(
    33761.901888296532,
    41398.852673195812,
    20294.441324124517,
    47050.795277378566,
    37831.274812427437,
    23805.972493309491,
    33269.521183912802,
    33269.521183912802,
    28200.836350389556,
    31202.006466498311,
    39832.479534257684,
    40498.81471219247,
    37520.657956701121,
    39367.632050026077,
    31940.938533600132,
    35138.695956433759,
    23714.90007790297,
    37738.26146697998,
    18225.064381868884,
    32660.622931255686,
    33014.789652200125,
    33020.856389657216,
    29847.28903130539,
    32306.847102126205,
    29787.725791887631,
    41031.006804096076,
    17058.906002364467,
    20975.410045154065,
    28801.097154515028,
    17345.210245339775,
    31399.932810423987,
    27171.376328810595,
    27680.46775989352,
    32778.506343622168,
    27680.46775989352,
    26802.117753919916,
    38865.893684203809,
    41470.569828189866,
    32086.363904946618,
    46331.372348461111,
    40533.459785764193,
    40245.463644812146,
    43552.955075216421,
    37635.685119755828,
    45113.805772148073,
    43574.699638304461,
    27017.007231062187,
    37737.294042691181,
    31280.724526515904,
    29684.448801554132,
    33280.457996465833,
    28706.117901219197,
    31280.724526515904,
    31924.681314365367,
    33333.537990300363,
    38672.792250745493,
    35346.844313648005,
    30866.827959939506,
    36467.337428263447,
    36244.364282188923,
    33980.834119373045,
    43666.324916988757,
    33968.280103946425,
    34997.501976300031,
    44930.242457986766,
    37737.294042691181,
    30610.116590385853,
    35898.30952283174,
    30713.378470011321,
    41964.064443456155,
    33803.610742146644,
    45337.673242382647,
    39937.914696019987,
    41838.104239785433,
    38150.061444343417,
    40104.724304045914,
    27871.214944358657,
    28369.155334313702,
    36035.219593904985,
    40189.166239912243,
    34529.989832063497,
    29482.013261411845,
    38623.584711271789,
    39129.600743199364,
    23844.103334863943,
    21358.037063633074,
    20955.540935780187,
    18351.438204667935,
    31303.528159702812,
    35533.752378073972,
    39334.26362616093,
    35342.436942704167,
    38283.688207204948,
    33437.021456288174,
    38639.23099651275,
    33870.640153459681,
    17359.041183041154,
    42679.617257571343,
    17718.645235463762,
    39381.385369783413,
    34521.026973957531,
    41400.430472176231,
    38144.85861434463,
    38249.806506617373,
    41400.430472176231,
    23578.716000934812,
    32877.164767615832,
    35411.626242926373,
    45179.955523152879,
    36857.900075876707,
    35686.462847723087,
    35686.462847723087,
    38597.935561185019,
    36850.370222252866,
    33820.549640786674,
    26133.753976839111,
    39972.98444509464,
    35535.491867163262,
    29325.729136202612,
    33148.843807430239,
    35278.858616282232,
    29008.130772101464,
    27605.807063044736,
    28091.381993220512,
    30633.20710686607,
    30536.393487826677,
    31316.747951567322,
    30633.20710686607,
    31300.96362253989,
    31300.96362253989,
    16951.029217865584,
    33770.51897116589,
    31300.96362253989,
    35669.461990702286,
    26871.561811127216,
    31300.96362253989,
    31538.744374715123,
    25530.33087812275,
    31300.96362253989,
    31312.160189492188,
    31644.261399543084,
    31312.160189492188,
    31300.96362253989,
    31300.96362253989,
    35669.461990702286,
    30712.621801222103,
    31312.160189492188,
    30651.04056501768,
    25530.33087812275,
    31300.96362253989,
    28000.549804436374,
    31312.160189492188,
    31300.96362253989,
    31312.160189492188,
    31300.96362253989,
    31300.96362253989,
    31081.746711899181
)

for ($i = 0; $i -lt $bValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $bValues[$i]
}
